# new button - refresh tickets / updating / v.12.7
# Adds the newly-collected ticket rows (71-77) coming from the
# "refresh tickets" run on 2024-05-20.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRows = @(
    @{ Row = 71; A = "2024-05-20"; B = "11:20:26"; C = "Etiquetadora";          D = "-"; E = "-"; F = "-"; G = "-"; H = "11:20:32"; I = "0:00:06" },
    @{ Row = 72; A = "2024-05-20"; B = "11:28:43"; C = "Fallo en elevador";     D = "-"; E = "-"; F = "-"; G = "-"; H = $null;      I = $null },
    @{ Row = 73; A = "2024-05-20"; B = "11:28:49"; C = "Fallo en paletizador";  D = "-"; E = "-"; F = "-"; G = "-"; H = "11:28:53"; I = "0:00:04" },
    @{ Row = 74; A = "2024-05-20"; B = "11:29:11"; C = "Ascensor no sube";      D = "-"; E = "-"; F = "-"; G = "-"; H = $null;      I = $null },
    @{ Row = 75; A = "2024-05-20"; B = "11:29:51"; C = "Fallo en paletizador";  D = "-"; E = "-"; F = "-"; G = "-"; H = "11:29:53"; I = "0:00:02" },
    @{ Row = 76; A = "2024-05-20"; B = "11:29:56"; C = "No atornilla clips";    D = "-"; E = "-"; F = "-"; G = "-"; H = "11:29:59"; I = "0:00:03" },
    @{ Row = 77; A = "2024-05-20"; B = "11:30:00"; C = "No pone tornillo";      D = "-"; E = "-"; F = "-"; G = "-"; H = $null;      I = $null }
)

foreach ($r in $newRows) {
    $row = $r.Row
    # Force column A to stay plain text (matches the rest of the sheet,
    # which stores the date as a literal string) instead of letting Excel
    # auto-convert the "yyyy-mm-dd" text into a date serial number. The
    # format is cleared again afterwards so the cell keeps the sheet's
    # default (General) style, same as every other row.
    $ws.Cells.Item($row, 1).NumberFormat = "@"
    $ws.Cells.Item($row, 1).Value = $r.A
    $ws.Cells.Item($row, 1).ClearFormats()
    $ws.Cells.Item($row, 2).Value = $r.B
    $ws.Cells.Item($row, 3).Value = $r.C
    $ws.Cells.Item($row, 4).Value = $r.D
    $ws.Cells.Item($row, 5).Value = $r.E
    $ws.Cells.Item($row, 6).Value = $r.F
    $ws.Cells.Item($row, 7).Value = $r.G
    if ($r.H) { $ws.Cells.Item($row, 8).Value = $r.H }
    if ($r.I) { $ws.Cells.Item($row, 9).Value = $r.I }
}
